$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price column (D) cells that look like plain numbers to remain text,
# matching the source data which stores prices as text strings (e.g. thousands
# separated with dots like "41.962.23").
$ws.Range("D2").Value = "41.962.23"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "2.209.23"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.96"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.98"
$ws.Range("E7").Value = "  -1.76%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.605"
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.87"
$ws.Range("E10").Value = "  +2.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0950"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.11"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("D14").Value = "2.542.85"
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.18"
$ws.Range("E15").Value = "  -2.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.838"
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").Value = "2.215.55"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "41.856.68"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000108"
$ws.Range("E19").Value = "  +9.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.68"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.12"
$ws.Range("E21").Value = "  -0.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.29"
$ws.Range("E22").Value = "  +18.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.02"
$ws.Range("E23").Value = "  -1.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("E24").Value = "  -7.60%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.44"
$ws.Range("E26").Value = "  +1.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.59"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.26"
$ws.Range("E28").Value = "  -2.59%  "
$ws.Range("E29").Value = "  +3.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.16"
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.51"
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.55"
$ws.Range("E32").Value = "  +7.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0789"
$ws.Range("E33").Value = "  -3.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.124"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "28.81"
$ws.Range("E35").Value = "  -4.97%  "
$ws.Range("E36").Value = "  -7.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.25"
$ws.Range("E37").Value = "  -5.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0300"
$ws.Range("E38").Value = "  -1.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "12.81"
$ws.Range("E39").Value = "  -6.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "65.44"
$ws.Range("E40").Value = "  +5.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.11"
$ws.Range("E41").Value = "  -3.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.61"
$ws.Range("E42").Value = "  -3.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.198"
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.68"
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.82"
$ws.Range("E45").Value = "  -3.46%  "
$ws.Range("E46").Value = "  -1.88%  "
$ws.Range("E47").Value = "  +4.82%  "
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.15"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").Value = "2.416.16"
$ws.Range("E51").Value = "  -1.44%  "
